# Update TPM-derived NATMI metrics for Adam9-Itgav LR pairs
# (commit: "update scripts wuth new tpm")
#
# This recomputation changes the ligand/receptor average & total expression
# values (and their derived specificity + edge weight columns) for rows 2-17
# on the active worksheet, reflecting new TPM inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.573375
$ws.Range("H2").Value = 28.720125
$ws.Range("I2").Value = 0.1037691388643484
$ws.Range("J2").Value = 0.1037691388643484
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 84.44336466475001
$ws.Range("R2").Value = 759.99028198275
$ws.Range("S2").Value = 0.006656942446394396
$ws.Range("T2").Value = 0.006656942446394396
$ws.Range("G3").Value = 9.573375
$ws.Range("H3").Value = 28.720125
$ws.Range("I3").Value = 0.1037691388643484
$ws.Range("J3").Value = 0.1037691388643484
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 523.774041380625
$ws.Range("R3").Value = 4713.966372425624
$ws.Range("S3").Value = 0.04129079486859573
$ws.Range("T3").Value = 0.04129079486859574
$ws.Range("G4").Value = 9.573375
$ws.Range("H4").Value = 28.720125
$ws.Range("I4").Value = 0.1037691388643484
$ws.Range("J4").Value = 0.1037691388643484
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 209.73506315125
$ws.Range("R4").Value = 1887.61556836125
$ws.Range("S4").Value = 0.01653409062904847
$ws.Range("T4").Value = 0.01653409062904847
$ws.Range("G5").Value = 9.573375
$ws.Range("H5").Value = 28.720125
$ws.Range("I5").Value = 0.1037691388643484
$ws.Range("J5").Value = 0.1037691388643484
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 498.359832529125
$ws.Range("R5").Value = 4485.238492762124
$ws.Range("S5").Value = 0.03928731092030979
$ws.Range("T5").Value = 0.03928731092030979
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("I6").Value = 0.5196887643218222
$ws.Range("J6").Value = 0.5196887643218222
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 422.9028815124736
$ws.Range("R6").Value = 3806.125933612263
$ws.Range("S6").Value = 0.03333879640892706
$ws.Range("T6").Value = 0.03333879640892706
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("I7").Value = 0.5196887643218222
$ws.Range("J7").Value = 0.5196887643218222
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("Q7").Value = 2623.125597146712
$ws.Range("R7").Value = 23608.13037432041
$ws.Range("S7").Value = 0.2067894404633893
$ws.Range("T7").Value = 0.2067894404633893
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("I8").Value = 0.5196887643218222
$ws.Range("J8").Value = 0.5196887643218222
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 1050.379303489432
$ws.Range("R8").Value = 9453.413731404891
$ws.Range("S8").Value = 0.08280478398715267
$ws.Range("T8").Value = 0.08280478398715267
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("I9").Value = 0.5196887643218222
$ws.Range("J9").Value = 0.5196887643218222
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 2495.8480757295
$ws.Range("R9").Value = 22462.6326815655
$ws.Range("S9").Value = 0.1967557434623531
$ws.Range("T9").Value = 0.1967557434623531
$ws.Range("G10").Value = 11.32006633333333
$ws.Range("H10").Value = 33.960199
$ws.Range("I10").Value = 0.122702133291269
$ws.Range("J10").Value = 0.122702133291269
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 99.85031291627313
$ws.Range("R10").Value = 898.6528162464581
$ws.Range("S10").Value = 0.007871521806088955
$ws.Range("T10").Value = 0.007871521806088955
$ws.Range("G11").Value = 11.32006633333333
$ws.Range("H11").Value = 33.960199
$ws.Range("I11").Value = 0.122702133291269
$ws.Range("J11").Value = 0.122702133291269
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 619.3382053984884
$ws.Range("R11").Value = 5574.043848586395
$ws.Range("S11").Value = 0.04882442575043423
$ws.Range("T11").Value = 0.04882442575043424
$ws.Range("G12").Value = 11.32006633333333
$ws.Range("H12").Value = 33.960199
$ws.Range("I12").Value = 0.122702133291269
$ws.Range("J12").Value = 0.122702133291269
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 248.0018621748345
$ws.Range("R12").Value = 2232.01675957351
$ws.Range("S12").Value = 0.01955078566150117
$ws.Range("T12").Value = 0.01955078566150117
$ws.Range("G13").Value = 11.32006633333333
$ws.Range("H13").Value = 33.960199
$ws.Range("I13").Value = 0.122702133291269
$ws.Range("J13").Value = 0.122702133291269
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 589.2870969849804
$ws.Range("R13").Value = 5303.583872864823
$ws.Range("S13").Value = 0.04645540007324459
$ws.Range("T13").Value = 0.0464554000732446
$ws.Range("G14").Value = 23.41838033333333
$ws.Range("H14").Value = 70.25514099999999
$ws.Range("I14").Value = 0.2538399635225604
$ws.Range("J14").Value = 0.2538399635225604
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 206.5652740382024
$ws.Range("R14").Value = 1859.087466343822
$ws.Range("S14").Value = 0.01628420594270823
$ws.Range("T14").Value = 0.01628420594270823
$ws.Range("G15").Value = 23.41838033333333
$ws.Range("H15").Value = 70.25514099999999
$ws.Range("I15").Value = 0.2538399635225604
$ws.Range("J15").Value = 0.2538399635225604
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 1281.255535250478
$ws.Range("R15").Value = 11531.2998172543
$ws.Range("S15").Value = 0.1010055010378705
$ws.Range("T15").Value = 0.1010055010378705
$ws.Range("G16").Value = 23.41838033333333
$ws.Range("H16").Value = 70.25514099999999
$ws.Range("I16").Value = 0.2538399635225604
$ws.Range("J16").Value = 0.2538399635225604
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 513.0537013447878
$ws.Range("R16").Value = 4617.48331210309
$ws.Range("S16").Value = 0.04044567593109636
$ws.Range("T16").Value = 0.04044567593109636
$ws.Range("G17").Value = 23.41838033333333
$ws.Range("H17").Value = 70.25514099999999
$ws.Range("I17").Value = 0.2538399635225604
$ws.Range("J17").Value = 0.2538399635225604
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 1219.087323020706
$ws.Range("R17").Value = 10971.78590718635
$ws.Range("S17").Value = 0.09610458061088537
$ws.Range("T17").Value = 0.09610458061088538
